# Edit script: add ownTeam / oppTeam columns and refresh the match-by-match
# data for the Marcus Stoinis (Delhi Capitals) IPL sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column D ("batsman") to make
# room for "ownTeam" and "oppTeam". This shifts batsman..sr from D:I to F:K.
$ws.Range("D:E").Insert()

# Header row
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Data rows (row, venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr)
$data = @(
    @(2,  " Dubai (DSC)",  " October 27 2020",   "Sunrisers won by 88 runs",                            "Delhi Capitals", "Sunrisers Hyderabad",         "Marcus Stoinis ", "5",  "6",  "1", "0", "83.33"),
    @(3,  " Dubai (DSC)",  " November 05 2020",  "Mumbai won by 57 runs",                                "Delhi Capitals", "Mumbai Indians",              "Marcus Stoinis ", "65", "46", "6", "3", "141.30"),
    @(4,  " Dubai (DSC)",  " November 10 2020",  "Mumbai won by 5 wickets (with 8 balls remaining)",     "Delhi Capitals", "Mumbai Indians",              "Marcus Stoinis ", "0",  "1",  "0", "0", "0.00"),
    @(5,  " Dubai (DSC)",  " October 14 2020",   "Capitals won by 13 runs",                              "Delhi Capitals", "Rajasthan Royals",            "Marcus Stoinis ", "18", "19", "1", "0", "94.73"),
    @(6,  " Dubai (DSC)",  " October 31 2020",   "Mumbai won by 9 wickets (with 34 balls remaining)",    "Delhi Capitals", "Mumbai Indians",              "Marcus Stoinis ", "2",  "3",  "0", "0", "66.66"),
    @(7,  " Dubai (DSC)",  " September 25 2020", "Capitals won by 44 runs",                              "Delhi Capitals", "Chennai Super Kings",         "Marcus Stoinis ", "5",  "3",  "1", "0", "166.66"),
    @(8,  " Sharjah",      " October 17 2020",   "Capitals won by 5 wickets (with 1 ball remaining)",    "Delhi Capitals", "Chennai Super Kings",         "Marcus Stoinis ", "24", "14", "1", "2", "171.42"),
    @(9,  " Sharjah",      " October 09 2020",   "Capitals won by 46 runs",                              "Delhi Capitals", "Rajasthan Royals",            "Marcus Stoinis ", "39", "30", "0", "4", "130.00"),
    @(10, " Dubai (DSC)",  " October 05 2020",   "Capitals won by 59 runs",                              "Delhi Capitals", "Royal Challengers Bangalore", "Marcus Stoinis ", "53", "26", "6", "2", "203.84"),
    @(11, " Dubai (DSC)",  " October 20 2020",   "Kings XI won by 5 wickets (with 6 balls remaining)",   "Delhi Capitals", "Kings XI Punjab",             "Marcus Stoinis ", "9",  "10", "0", "0", "90.00"),
    @(12, " Dubai (DSC)",  " September 20 2020", "Match tied (Capitals won the one-over eliminator)",    "Delhi Capitals", "Kings XI Punjab",             "Marcus Stoinis ", "53", "21", "7", "3", "252.38"),
    @(13, " Sharjah",      " October 03 2020",   "Capitals won by 18 runs",                              "Delhi Capitals", "Kolkata Knight Riders",       "Marcus Stoinis ", "1",  "3",  "0", "0", "33.33"),
    @(14, " Abu Dhabi",    " November 08 2020",  "Capitals won by 17 runs",                              "Delhi Capitals", "Sunrisers Hyderabad",         "Marcus Stoinis ", "38", "27", "5", "1", "140.74"),
    @(15, " Abu Dhabi",    " November 02 2020",  "Capitals won by 6 wickets (with 6 balls remaining)",   "Delhi Capitals", "Royal Challengers Bangalore", "Marcus Stoinis ", "10", "5",  "0", "1", "200.00"),
    @(16, " Abu Dhabi",    " October 24 2020",   "KKR won by 59 runs",                                   "Delhi Capitals", "Kolkata Knight Riders",       "Marcus Stoinis ", "6",  "6",  "0", "0", "100.00"),
    @(17, " Abu Dhabi",    " October 11 2020",   "Mumbai won by 5 wickets (with 2 balls remaining)",     "Delhi Capitals", "Mumbai Indians",              "Marcus Stoinis ", "13", "8",  "2", "0", "162.50"),
    @(18, " Abu Dhabi",    " September 29 2020", "Sunrisers won by 15 runs",                             "Delhi Capitals", "Sunrisers Hyderabad",         "Marcus Stoinis ", "11", "9",  "1", "0", "122.22")
)

# Columns G:K (totalRuns, totalBalls, total4s, total6s, sr) hold numeric
# looking values that must stay text (matches the workbook's existing
# "number stored as text" convention) -- format as Text first so Excel
# doesn't silently coerce them (and lose precision on the sr column).
$ws.Range("G2:K18").NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value  = $row[1]   # A venue
    $ws.Cells.Item($r, 2).Value  = $row[2]   # B date
    $ws.Cells.Item($r, 3).Value  = $row[3]   # C result
    $ws.Cells.Item($r, 4).Value  = $row[4]   # D ownTeam
    $ws.Cells.Item($r, 5).Value  = $row[5]   # E oppTeam
    $ws.Cells.Item($r, 6).Value  = $row[6]   # F batsman
    $ws.Cells.Item($r, 7).Value  = $row[7]   # G totalRuns
    $ws.Cells.Item($r, 8).Value  = $row[8]   # H totalBalls
    $ws.Cells.Item($r, 9).Value  = $row[9]   # I total4s
    $ws.Cells.Item($r, 10).Value = $row[10]  # J total6s
    $ws.Cells.Item($r, 11).Value = $row[11]  # K sr
}
